$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new date blocks (Status/Time) for 09-04-2025 and 10-04-2025
# mirroring the existing 07-04-2025 Status/Time columns (D:E -> F:G, H:I)

# Copy formatting from the existing header cells (D1:E1) to the new header
# cells so the style (bold, border, centered) matches.
$ws.Range("D1:E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("F1").Value = "09-04-2025 Status"
$ws.Range("G1").Value = "09-04-2025 Time"
$ws.Range("H1").Value = "10-04-2025 Status"
$ws.Range("I1").Value = "10-04-2025 Time"

for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 6).Value = "A"          # F - 09-04-2025 Status
    $ws.Cells.Item($r, 7).Value = "00:00:00"   # G - 09-04-2025 Time
    $ws.Cells.Item($r, 8).Value = "A"          # H - 10-04-2025 Status
    $ws.Cells.Item($r, 9).Value = "00:00:00"   # I - 10-04-2025 Time
}
